$d = $word.ActiveDocument

# The sentence currently reads:
#   "...a burst of ultrasound (in audible sounds) to the obstacle..."
# and must become:
#   "...a burst of ultrasound to the obstacle..."
# i.e. remove the parenthetical "(in audible sounds) " (including the
# trailing space so the remaining words stay correctly spaced).
$d.Content.Find.Execute("(in audible sounds) ", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
